# Updating price policy code
# The "Electric_boiler" technology row is removed from the per-technology cost
# sheets (its label row moves off the bottom of each of those tables) and the
# remaining technology costs are recalculated.

$wb = $excel.ActiveWorkbook

# --- Sheets 1-4: Operating / Maintenance / Capital / Total cost per technology ---
# Each of these sheets has 7 rows (Electric_boiler, Gas_CHP, Gas_boiler, Grid,
# Heat_pump, Solar_PV, Solar_thermal). The "Electric_boiler" technology is
# removed from the model entirely, so its row is deleted (the remaining rows
# shift up, leaving a 6-row table), and new recalculated cost values are
# written in for the remaining technologies.

$sheetNames = @(
    "Operating_cost_per_technology",
    "Maintenance_cost_per_technology",
    "Capital_cost_per_technology",
    "Total_cost_per_technology"
)

$newValues = @{
    "Operating_cost_per_technology"   = @(36192.078682714076, 104121.56698515784, 0, 0, 0, 0)
    "Maintenance_cost_per_technology" = @(6916.3062362667724, 10874.919218450003, 0, 0, 31336.456223487592, 0)
    "Capital_cost_per_technology"     = @(6474.830691193004, 35345.246760457267, 0, 0, 76114.007979724091, 0)
    "Total_cost_per_technology"       = @(49583.215610173851, 150341.7329640651, 0, 0, 107450.46420321168, 0)
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Drop the "Electric_boiler" row (row 1); this shifts every remaining
    # technology row up by one, leaving a 6-row table.
    $ws.Rows.Item(1).Delete()

    $vals = $newValues[$name]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($i + 1, 2).Value = $vals[$i]
    }
}

# --- Sheet 5 & 6: Operating_cost_grid / Total_cost_grid ---
$wb.Worksheets.Item("Operating_cost_grid").Range("A1").Value = 143167.96235087828
$wb.Worksheets.Item("Total_cost_grid").Range("A1").Value = 143167.96235087828

# --- Sheet 7 & 8: Capital_cost_per_storage / Total_cost_per_storage ---
$wb.Worksheets.Item("Capital_cost_per_storage").Range("B2").Value = 3656.9132202519372
$wb.Worksheets.Item("Total_cost_per_storage").Range("B2").Value = 3656.9132202519372

# --- Sheet 9: Income_via_exports --- (unchanged)
